$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H: build the SQL "('<name>')," literal list next to the existing
# DECLARE / CONCATENATE helper columns. H2 is a standalone formula; H3:H48
# is entered as one range assignment so Excel collapses it into a shared
# formula group (mirrors columns D/E already on the sheet).
$ws.Range("H2").Formula = '="("&"''"&A2&"''"&")"&","'
$ws.Range("H3:H48").Formula = '="("&"''"&A3&"''"&")"&","'

# Move the selection / scroll position onto the new column, like the
# author did after adding it.
$ws.Range("H2:H48").Select()
$excel.ActiveWindow.ScrollRow = 23
